# Planteamiento.docx - "Definición de la clase proyectil"
#
# 1) Fix a couple of existing paragraphs (typo fix / add signature).
# 2) Prefix the "Proyectil_Fisicas" attribute list with "Float " and fix a
#    typo ("Velocicad" -> "Velocidad").
# 3) Append the full "Proyectil_Fisicas" method list, plus the whole new
#    "Proyectil_Graph" class (attributes + methods).

$d = $word.ActiveDocument

function Set-ParaText([int]$index, [string]$newText) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null   # keep the paragraph mark out of the replace
    $r.Text = $newText
}

function Add-ListPara([string]$text, [int]$level) {
    # Clone formatting/list level from the current last paragraph, then
    # retarget the level and set the text.
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Last
    $newPara.Range.ListFormat.ListLevelNumber = $level
    $r = $newPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

# --- 1) "cañón" class: fix constructor typo/signature -----------------
Set-ParaText 21 "Constructor (float posX, float posY)"

# --- 2) "Proyectil_Fisicas" attributes: add "Float " / fix typo --------
Set-ParaText 27 "Float Posición en X"
Set-ParaText 28 "Float Posición en Y"
Set-ParaText 29 "Float Ángulo"
Set-ParaText 30 "Float Velocidad inicial "
Set-ParaText 31 "Float Velocidad en X"
Set-ParaText 32 "Float Velocidad en Y"
Set-ParaText 33 "Float Radio"
Set-ParaText 34 "Float Aceleración en X"
Set-ParaText 35 "Float Aceleración en Y"
Set-ParaText 36 "Float Delta de tiempo"

# --- 3) New attributes + "Métodos" for Proyectil_Fisicas ----------------
Add-ListPara "Float Rango de daño" 4
Add-ListPara "Int Identificador, para variar entre proyectil ofensivo y defensivo" 4
Add-ListPara "Métodos " 3
Add-ListPara "Constructor (float Px, float Py, float Angulo, float velocidad inicial, float radio)" 4
Add-ListPara "Destructor " 4
Add-ListPara "Actualizar Posición" 4
Add-ListPara "Calcular velocidad" 4
Add-ListPara " Métodos get y set necesarios" 4

# --- New class: Proyectil_Graph -----------------------------------------
Add-ListPara "Proyectil_Graph" 2
Add-ListPara "Atributos" 3
Add-ListPara "Instancia de la clase Proyectil_fisicas" 4
Add-ListPara "Int Id, para variar entre proyectil ofensivo o defensivo" 4
Add-ListPara "Métodos " 3
Add-ListPara "Constructor (float Px, float Py, float Vel_In, float angle, float range)" 4
